$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prevent Excel from auto-converting the Date column (B) text into a date serial number
$ws.Range("B2:B10").NumberFormat = "@"

# Row 2: FK Vrsac vs Macva Sabac
$ws.Cells.Item(2, 1).Value2 = "Serbian First League"
$ws.Cells.Item(2, 2).Value2 = "2025-11-17"
$ws.Cells.Item(2, 3).Value2 = "09:00:00"
$ws.Cells.Item(2, 4).Value2 = "FK Vrsac"
$ws.Cells.Item(2, 5).Value2 = "Macva Sabac"
$ws.Cells.Item(2, 6).Value2 = 1.02
$ws.Cells.Item(2, 7).Value2 = 1000
$ws.Cells.Item(2, 8).Value2 = 1.02
$ws.Cells.Item(2, 9).Value2 = 1000
$ws.Cells.Item(2, 10).Value2 = 1.02
$ws.Cells.Item(2, 11).Value2 = 1000
$ws.Cells.Item(2, 12).Value2 = 1.01
$ws.Cells.Item(2, 13).Value2 = 1.01
$ws.Cells.Item(2, 14).Value2 = 1.1
$ws.Cells.Item(2, 15).Value2 = 1.01
$ws.Cells.Item(2, 16).Value2 = 1.1
$ws.Cells.Item(2, 17).Value2 = 1.01
$ws.Cells.Item(2, 18).Value2 = 1.09
$ws.Cells.Item(2, 19).Value2 = 1.02
$ws.Cells.Item(2, 20).Value2 = 1.01
$ws.Cells.Item(2, 21).Value2 = 1.01
$ws.Cells.Item(2, 22).Value2 = 1.01
$ws.Cells.Item(2, 23).Value2 = 1.01
$ws.Cells.Item(2, 24).Value2 = 1000
$ws.Cells.Item(2, 25).Value2 = 1000
$ws.Cells.Item(2, 26).Value2 = 1000
$ws.Cells.Item(2, 27).Value2 = 1000
$ws.Cells.Item(2, 28).Value2 = 1000
$ws.Cells.Item(2, 29).Value2 = 1000
$ws.Cells.Item(2, 30).Value2 = 1000
$ws.Cells.Item(2, 31).Value2 = 1000
$ws.Cells.Item(2, 32).Value2 = 1000
$ws.Cells.Item(2, 33).Value2 = 1000
$ws.Cells.Item(2, 34).Value2 = 1000
$ws.Cells.Item(2, 35).Value2 = 1000
$ws.Cells.Item(2, 36).Value2 = 1000
$ws.Cells.Item(2, 37).Value2 = 1000
$ws.Cells.Item(2, 38).Value2 = 1000
$ws.Cells.Item(2, 39).Value2 = 1000
$ws.Cells.Item(2, 40).Value2 = 1000
$ws.Cells.Item(2, 41).Value2 = 1000

# Row 3: Tekstilac Odzaci vs FK Kabel
$ws.Cells.Item(3, 1).Value2 = "Serbian First League"
$ws.Cells.Item(3, 2).Value2 = "2025-11-17"
$ws.Cells.Item(3, 3).Value2 = "09:00:00"
$ws.Cells.Item(3, 4).Value2 = "Tekstilac Odzaci"
$ws.Cells.Item(3, 5).Value2 = "FK Kabel"
$ws.Cells.Item(3, 6).Value2 = 1.02
$ws.Cells.Item(3, 7).Value2 = 1000
$ws.Cells.Item(3, 8).Value2 = 1.02
$ws.Cells.Item(3, 9).Value2 = 1000
$ws.Cells.Item(3, 10).Value2 = 1.02
$ws.Cells.Item(3, 11).Value2 = 1000
$ws.Cells.Item(3, 12).Value2 = 1.01
$ws.Cells.Item(3, 13).Value2 = 1.01
$ws.Cells.Item(3, 14).Value2 = 1.1
$ws.Cells.Item(3, 15).Value2 = 1.01
$ws.Cells.Item(3, 16).Value2 = 1.1
$ws.Cells.Item(3, 17).Value2 = 1.01
$ws.Cells.Item(3, 18).Value2 = 1.09
$ws.Cells.Item(3, 19).Value2 = 1.02
$ws.Cells.Item(3, 20).Value2 = 1.01
$ws.Cells.Item(3, 21).Value2 = 1.01
$ws.Cells.Item(3, 22).Value2 = 1.01
$ws.Cells.Item(3, 23).Value2 = 1.01
$ws.Cells.Item(3, 24).Value2 = 1000
$ws.Cells.Item(3, 25).Value2 = 1000
$ws.Cells.Item(3, 26).Value2 = 1000
$ws.Cells.Item(3, 27).Value2 = 1000
$ws.Cells.Item(3, 28).Value2 = 1000
$ws.Cells.Item(3, 29).Value2 = 1000
$ws.Cells.Item(3, 30).Value2 = 1000
$ws.Cells.Item(3, 31).Value2 = 1000
$ws.Cells.Item(3, 32).Value2 = 1000
$ws.Cells.Item(3, 33).Value2 = 1000
$ws.Cells.Item(3, 34).Value2 = 1000
$ws.Cells.Item(3, 35).Value2 = 1000
$ws.Cells.Item(3, 36).Value2 = 1000
$ws.Cells.Item(3, 37).Value2 = 1000
$ws.Cells.Item(3, 38).Value2 = 1000
$ws.Cells.Item(3, 39).Value2 = 1000
$ws.Cells.Item(3, 40).Value2 = 1000
$ws.Cells.Item(3, 41).Value2 = 1000

# Row 4: ES Mostaganem vs USM Alger
$ws.Cells.Item(4, 1).Value2 = "Algerian Ligue 1"
$ws.Cells.Item(4, 2).Value2 = "2025-11-17"
$ws.Cells.Item(4, 3).Value2 = "12:00:00"
$ws.Cells.Item(4, 4).Value2 = "ES Mostaganem"
$ws.Cells.Item(4, 5).Value2 = "USM Alger"
$ws.Cells.Item(4, 6).Value2 = 4.7
$ws.Cells.Item(4, 7).Value2 = 5.4
$ws.Cells.Item(4, 8).Value2 = 1.87
$ws.Cells.Item(4, 9).Value2 = 2.04
$ws.Cells.Item(4, 10).Value2 = 3.05
$ws.Cells.Item(4, 11).Value2 = 3.6
$ws.Cells.Item(4, 12).Value2 = 1.51
$ws.Cells.Item(4, 13).Value2 = 1.12
$ws.Cells.Item(4, 14).Value2 = 2.3
$ws.Cells.Item(4, 15).Value2 = 1.59
$ws.Cells.Item(4, 16).Value2 = 1.44
$ws.Cells.Item(4, 17).Value2 = 2.52
$ws.Cells.Item(4, 18).Value2 = 1.15
$ws.Cells.Item(4, 19).Value2 = 3.45
$ws.Cells.Item(4, 20).Value2 = 2.34
$ws.Cells.Item(4, 21).Value2 = 1.59
$ws.Cells.Item(4, 22).Value2 = 1.96
$ws.Cells.Item(4, 23).Value2 = 1.22
$ws.Cells.Item(4, 24).Value2 = 8.4
$ws.Cells.Item(4, 25).Value2 = 6.2
$ws.Cells.Item(4, 26).Value2 = 10.5
$ws.Cells.Item(4, 27).Value2 = 25
$ws.Cells.Item(4, 28).Value2 = 14.5
$ws.Cells.Item(4, 29).Value2 = 7.8
$ws.Cells.Item(4, 30).Value2 = 12.5
$ws.Cells.Item(4, 31).Value2 = 32
$ws.Cells.Item(4, 32).Value2 = 42
$ws.Cells.Item(4, 33).Value2 = 28
$ws.Cells.Item(4, 34).Value2 = 36
$ws.Cells.Item(4, 35).Value2 = 85
$ws.Cells.Item(4, 36).Value2 = 190
$ws.Cells.Item(4, 37).Value2 = 150
$ws.Cells.Item(4, 38).Value2 = 180
$ws.Cells.Item(4, 39).Value2 = 350
$ws.Cells.Item(4, 40).Value2 = 260
$ws.Cells.Item(4, 41).Value2 = 28

# Row 5: Foggia vs Cavese 1919
$ws.Cells.Item(5, 1).Value2 = "Italian Serie C"
$ws.Cells.Item(5, 2).Value2 = "2025-11-17"
$ws.Cells.Item(5, 3).Value2 = "16:30:00"
$ws.Cells.Item(5, 4).Value2 = "Foggia"
$ws.Cells.Item(5, 5).Value2 = "Cavese 1919"
$ws.Cells.Item(5, 6).Value2 = 3.6
$ws.Cells.Item(5, 7).Value2 = 4.6
$ws.Cells.Item(5, 8).Value2 = 2.1
$ws.Cells.Item(5, 9).Value2 = 2.44
$ws.Cells.Item(5, 10).Value2 = 2.94
$ws.Cells.Item(5, 11).Value2 = 3.45
$ws.Cells.Item(5, 12).Value2 = 1.45
$ws.Cells.Item(5, 13).Value2 = 1.1
$ws.Cells.Item(5, 14).Value2 = 2.58
$ws.Cells.Item(5, 15).Value2 = 1.47
$ws.Cells.Item(5, 16).Value2 = 1.54
$ws.Cells.Item(5, 17).Value2 = 2.42
$ws.Cells.Item(5, 18).Value2 = 1.2
$ws.Cells.Item(5, 19).Value2 = 4.5
$ws.Cells.Item(5, 20).Value2 = 2
$ws.Cells.Item(5, 21).Value2 = 1.76
$ws.Cells.Item(5, 22).Value2 = 1.7
$ws.Cells.Item(5, 23).Value2 = 1.28
$ws.Cells.Item(5, 24).Value2 = 1000
$ws.Cells.Item(5, 25).Value2 = 1000
$ws.Cells.Item(5, 26).Value2 = 1000
$ws.Cells.Item(5, 27).Value2 = 1000
$ws.Cells.Item(5, 28).Value2 = 1000
$ws.Cells.Item(5, 29).Value2 = 1000
$ws.Cells.Item(5, 30).Value2 = 1000
$ws.Cells.Item(5, 31).Value2 = 1000
$ws.Cells.Item(5, 32).Value2 = 1000
$ws.Cells.Item(5, 33).Value2 = 1000
$ws.Cells.Item(5, 34).Value2 = 1000
$ws.Cells.Item(5, 35).Value2 = 1000
$ws.Cells.Item(5, 36).Value2 = 1000
$ws.Cells.Item(5, 37).Value2 = 1000
$ws.Cells.Item(5, 38).Value2 = 1000
$ws.Cells.Item(5, 39).Value2 = 1000
$ws.Cells.Item(5, 40).Value2 = 1000
$ws.Cells.Item(5, 41).Value2 = 1000

# Row 6: Pergolettese vs Giana Erminio
$ws.Cells.Item(6, 1).Value2 = "Italian Serie C"
$ws.Cells.Item(6, 2).Value2 = "2025-11-17"
$ws.Cells.Item(6, 3).Value2 = "16:30:00"
$ws.Cells.Item(6, 4).Value2 = "Pergolettese"
$ws.Cells.Item(6, 5).Value2 = "Giana Erminio"
$ws.Cells.Item(6, 6).Value2 = 2.36
$ws.Cells.Item(6, 7).Value2 = 3.3
$ws.Cells.Item(6, 8).Value2 = 2.66
$ws.Cells.Item(6, 9).Value2 = 3.85
$ws.Cells.Item(6, 10).Value2 = 2.64
$ws.Cells.Item(6, 11).Value2 = 3.7
$ws.Cells.Item(6, 12).Value2 = 1.52
$ws.Cells.Item(6, 13).Value2 = 1.01
$ws.Cells.Item(6, 14).Value2 = 1.38
$ws.Cells.Item(6, 15).Value2 = 1.01
$ws.Cells.Item(6, 16).Value2 = 1.38
$ws.Cells.Item(6, 17).Value2 = 2.52
$ws.Cells.Item(6, 18).Value2 = 1.18
$ws.Cells.Item(6, 19).Value2 = 2.54
$ws.Cells.Item(6, 20).Value2 = 1.01
$ws.Cells.Item(6, 21).Value2 = 1.01
$ws.Cells.Item(6, 22).Value2 = 1.35
$ws.Cells.Item(6, 23).Value2 = 1.44
$ws.Cells.Item(6, 24).Value2 = 1000
$ws.Cells.Item(6, 25).Value2 = 1000
$ws.Cells.Item(6, 26).Value2 = 1000
$ws.Cells.Item(6, 27).Value2 = 1000
$ws.Cells.Item(6, 28).Value2 = 1000
$ws.Cells.Item(6, 29).Value2 = 1000
$ws.Cells.Item(6, 30).Value2 = 1000
$ws.Cells.Item(6, 31).Value2 = 1000
$ws.Cells.Item(6, 32).Value2 = 1000
$ws.Cells.Item(6, 33).Value2 = 1000
$ws.Cells.Item(6, 34).Value2 = 1000
$ws.Cells.Item(6, 35).Value2 = 1000
$ws.Cells.Item(6, 36).Value2 = 1000
$ws.Cells.Item(6, 37).Value2 = 1000
$ws.Cells.Item(6, 38).Value2 = 1000
$ws.Cells.Item(6, 39).Value2 = 1000
$ws.Cells.Item(6, 40).Value2 = 1000
$ws.Cells.Item(6, 41).Value2 = 1000

# Row 7: Barracas Central vs Huracan
$ws.Cells.Item(7, 1).Value2 = "Argentinian Primera Division"
$ws.Cells.Item(7, 2).Value2 = "2025-11-17"
$ws.Cells.Item(7, 3).Value2 = "17:00:00"
$ws.Cells.Item(7, 4).Value2 = "Barracas Central"
$ws.Cells.Item(7, 5).Value2 = "Huracan"
$ws.Cells.Item(7, 6).Value2 = 3.3
$ws.Cells.Item(7, 7).Value2 = 3.4
$ws.Cells.Item(7, 8).Value2 = 2.82
$ws.Cells.Item(7, 9).Value2 = 2.92
$ws.Cells.Item(7, 10).Value2 = 2.82
$ws.Cells.Item(7, 11).Value2 = 2.84
$ws.Cells.Item(7, 12).Value2 = 1.01
$ws.Cells.Item(7, 13).Value2 = 1.19
$ws.Cells.Item(7, 14).Value2 = 2.14
$ws.Cells.Item(7, 15).Value2 = 1.82
$ws.Cells.Item(7, 16).Value2 = 1.34
$ws.Cells.Item(7, 17).Value2 = 3.6
$ws.Cells.Item(7, 18).Value2 = 1.11
$ws.Cells.Item(7, 19).Value2 = 8.8
$ws.Cells.Item(7, 20).Value2 = 2.64
$ws.Cells.Item(7, 21).Value2 = 1.52
$ws.Cells.Item(7, 22).Value2 = 1.52
$ws.Cells.Item(7, 23).Value2 = 1.41
$ws.Cells.Item(7, 24).Value2 = 5.8
$ws.Cells.Item(7, 25).Value2 = 7
$ws.Cells.Item(7, 26).Value2 = 15
$ws.Cells.Item(7, 27).Value2 = 150
$ws.Cells.Item(7, 28).Value2 = 7.6
$ws.Cells.Item(7, 29).Value2 = 7.6
$ws.Cells.Item(7, 30).Value2 = 19
$ws.Cells.Item(7, 31).Value2 = 60
$ws.Cells.Item(7, 32).Value2 = 19.5
$ws.Cells.Item(7, 33).Value2 = 24
$ws.Cells.Item(7, 34).Value2 = 36
$ws.Cells.Item(7, 35).Value2 = 130
$ws.Cells.Item(7, 36).Value2 = 80
$ws.Cells.Item(7, 37).Value2 = 75
$ws.Cells.Item(7, 38).Value2 = 140
$ws.Cells.Item(7, 39).Value2 = 370
$ws.Cells.Item(7, 40).Value2 = 130
$ws.Cells.Item(7, 41).Value2 = 1000

# Row 8: Belgrano vs Union Santa Fe
$ws.Cells.Item(8, 1).Value2 = "Argentinian Primera Division"
$ws.Cells.Item(8, 2).Value2 = "2025-11-17"
$ws.Cells.Item(8, 3).Value2 = "17:00:00"
$ws.Cells.Item(8, 4).Value2 = "Belgrano"
$ws.Cells.Item(8, 5).Value2 = "Union Santa Fe"
$ws.Cells.Item(8, 6).Value2 = 2.22
$ws.Cells.Item(8, 7).Value2 = 2.36
$ws.Cells.Item(8, 8).Value2 = 3.85
$ws.Cells.Item(8, 9).Value2 = 4.2
$ws.Cells.Item(8, 10).Value2 = 3.05
$ws.Cells.Item(8, 11).Value2 = 3.3
$ws.Cells.Item(8, 12).Value2 = 1.61
$ws.Cells.Item(8, 13).Value2 = 1.13
$ws.Cells.Item(8, 14).Value2 = 2.58
$ws.Cells.Item(8, 15).Value2 = 1.57
$ws.Cells.Item(8, 16).Value2 = 1.51
$ws.Cells.Item(8, 17).Value2 = 2.78
$ws.Cells.Item(8, 18).Value2 = 1.18
$ws.Cells.Item(8, 19).Value2 = 5.5
$ws.Cells.Item(8, 20).Value2 = 2.12
$ws.Cells.Item(8, 21).Value2 = 1.73
$ws.Cells.Item(8, 22).Value2 = 1.31
$ws.Cells.Item(8, 23).Value2 = 1.73
$ws.Cells.Item(8, 24).Value2 = 9.8
$ws.Cells.Item(8, 25).Value2 = 980
$ws.Cells.Item(8, 26).Value2 = 980
$ws.Cells.Item(8, 27).Value2 = 120
$ws.Cells.Item(8, 28).Value2 = 980
$ws.Cells.Item(8, 29).Value2 = 980
$ws.Cells.Item(8, 30).Value2 = 980
$ws.Cells.Item(8, 31).Value2 = 90
$ws.Cells.Item(8, 32).Value2 = 980
$ws.Cells.Item(8, 33).Value2 = 980
$ws.Cells.Item(8, 34).Value2 = 980
$ws.Cells.Item(8, 35).Value2 = 120
$ws.Cells.Item(8, 36).Value2 = 980
$ws.Cells.Item(8, 37).Value2 = 980
$ws.Cells.Item(8, 38).Value2 = 80
$ws.Cells.Item(8, 39).Value2 = 1000
$ws.Cells.Item(8, 40).Value2 = 980
$ws.Cells.Item(8, 41).Value2 = 1000

# Row 9: Defensa y Justicia vs Independiente Rivadavia
$ws.Cells.Item(9, 1).Value2 = "Argentinian Primera Division"
$ws.Cells.Item(9, 2).Value2 = "2025-11-17"
$ws.Cells.Item(9, 3).Value2 = "17:00:00"
$ws.Cells.Item(9, 4).Value2 = "Defensa y Justicia"
$ws.Cells.Item(9, 5).Value2 = "Independiente Rivadavia"
$ws.Cells.Item(9, 6).Value2 = 1.99
$ws.Cells.Item(9, 7).Value2 = 2.14
$ws.Cells.Item(9, 8).Value2 = 4.4
$ws.Cells.Item(9, 9).Value2 = 5.1
$ws.Cells.Item(9, 10).Value2 = 3.15
$ws.Cells.Item(9, 11).Value2 = 3.45
$ws.Cells.Item(9, 12).Value2 = 1.46
$ws.Cells.Item(9, 13).Value2 = 1.11
$ws.Cells.Item(9, 14).Value2 = 2.7
$ws.Cells.Item(9, 15).Value2 = 1.5
$ws.Cells.Item(9, 16).Value2 = 1.57
$ws.Cells.Item(9, 17).Value2 = 2.28
$ws.Cells.Item(9, 18).Value2 = 1.21
$ws.Cells.Item(9, 19).Value2 = 4.4
$ws.Cells.Item(9, 20).Value2 = 2.12
$ws.Cells.Item(9, 21).Value2 = 1.78
$ws.Cells.Item(9, 22).Value2 = 1.25
$ws.Cells.Item(9, 23).Value2 = 1.89
$ws.Cells.Item(9, 24).Value2 = 1000
$ws.Cells.Item(9, 25).Value2 = 980
$ws.Cells.Item(9, 26).Value2 = 1000
$ws.Cells.Item(9, 27).Value2 = 1000
$ws.Cells.Item(9, 28).Value2 = 8.2
$ws.Cells.Item(9, 29).Value2 = 980
$ws.Cells.Item(9, 30).Value2 = 980
$ws.Cells.Item(9, 31).Value2 = 1000
$ws.Cells.Item(9, 32).Value2 = 980
$ws.Cells.Item(9, 33).Value2 = 980
$ws.Cells.Item(9, 34).Value2 = 980
$ws.Cells.Item(9, 35).Value2 = 120
$ws.Cells.Item(9, 36).Value2 = 980
$ws.Cells.Item(9, 37).Value2 = 980
$ws.Cells.Item(9, 38).Value2 = 1000
$ws.Cells.Item(9, 39).Value2 = 1000
$ws.Cells.Item(9, 40).Value2 = 32
$ws.Cells.Item(9, 41).Value2 = 1000

# Row 10: CA Platense vs Gimnasia La Plata
$ws.Cells.Item(10, 1).Value2 = "Argentinian Primera Division"
$ws.Cells.Item(10, 2).Value2 = "2025-11-17"
$ws.Cells.Item(10, 3).Value2 = "19:30:00"
$ws.Cells.Item(10, 4).Value2 = "CA Platense"
$ws.Cells.Item(10, 5).Value2 = "Gimnasia La Plata"
$ws.Cells.Item(10, 6).Value2 = 2.42
$ws.Cells.Item(10, 7).Value2 = 2.5
$ws.Cells.Item(10, 8).Value2 = 3.35
$ws.Cells.Item(10, 9).Value2 = 3.45
$ws.Cells.Item(10, 10).Value2 = 3.2
$ws.Cells.Item(10, 11).Value2 = 3.4
$ws.Cells.Item(10, 12).Value2 = 1.01
$ws.Cells.Item(10, 13).Value2 = 1.12
$ws.Cells.Item(10, 14).Value2 = 2.52
$ws.Cells.Item(10, 15).Value2 = 1.6
$ws.Cells.Item(10, 16).Value2 = 1.49
$ws.Cells.Item(10, 17).Value2 = 2.72
$ws.Cells.Item(10, 18).Value2 = 1.17
$ws.Cells.Item(10, 19).Value2 = 6.2
$ws.Cells.Item(10, 20).Value2 = 2.26
$ws.Cells.Item(10, 21).Value2 = 1.69
$ws.Cells.Item(10, 22).Value2 = 1.4
$ws.Cells.Item(10, 23).Value2 = 1.66
$ws.Cells.Item(10, 24).Value2 = 8.6
$ws.Cells.Item(10, 25).Value2 = 9.2
$ws.Cells.Item(10, 26).Value2 = 21
$ws.Cells.Item(10, 27).Value2 = 75
$ws.Cells.Item(10, 28).Value2 = 7.4
$ws.Cells.Item(10, 29).Value2 = 7.8
$ws.Cells.Item(10, 30).Value2 = 17
$ws.Cells.Item(10, 31).Value2 = 65
$ws.Cells.Item(10, 32).Value2 = 13
$ws.Cells.Item(10, 33).Value2 = 16
$ws.Cells.Item(10, 34).Value2 = 27
$ws.Cells.Item(10, 35).Value2 = 110
$ws.Cells.Item(10, 36).Value2 = 60
$ws.Cells.Item(10, 37).Value2 = 70
$ws.Cells.Item(10, 38).Value2 = 75
$ws.Cells.Item(10, 39).Value2 = 240
$ws.Cells.Item(10, 40).Value2 = 44
$ws.Cells.Item(10, 41).Value2 = 95
